$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "--"
$ws.Range("B12").Value = "ffff1000DBD01E69CBD0E196"
$ws.Range("C12").Value = "Unknown"
